$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Remove the "Meta description: ..." paragraph that currently sits
#    right after the title (Heading1) paragraph.
# ------------------------------------------------------------------
$metaPara = $d.Paragraphs.Item(2)
if ($metaPara.Range.Text -like "Meta description*") {
    $metaPara.Range.Delete()
}

# ------------------------------------------------------------------
# 2. Insert a new paragraph - bold "Play Boom Pirates Slot - Free
#    Online Review" - right before the final ("For the feature
#    image...") paragraph.
# ------------------------------------------------------------------
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count)
$newParaRange = $lastPara.Range.InsertParagraphBefore()
$insertPoint = $d.Paragraphs.Item($count).Range

$xml = @"
<?xml version="1.0" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>
<w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Boom Pirates Slot - Free Online Review</w:t></w:r></w:p>
</w:body></w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@

$insertPoint.InsertXML($xml)

# ------------------------------------------------------------------
# 3. Swap out the text of the (now last) "For the feature image..."
#    paragraph for the new marketing blurb, keeping its italic run
#    formatting intact.
# ------------------------------------------------------------------
$oldText = "For the feature image of Boom Pirates, we want to highlight the main character of the game, Captain Mary. In a cartoon style, create an image of a happy Maya warrior wearing glasses, with a pirate hat and a peg leg, holding a treasure map and a sword. In the background, include a pirate ship sailing the high seas with the Jolly Roger flag flying high. The overall theme should be vibrant, colorful, and action-packed to capture the exciting and adventurous feeling of the game."
$newText = "Join Captain Mary and her band of pirates in the Boom Pirates slot - a thrilling game with innovative mechanics, bonuses, and potential payouts of 3,800 times the bet. Play now for free."

$d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
